# Update cryptocurrency price (column D) and 1h volume change (column E)
# values to reflect the latest scrape, as produced by the scheduled
# GitHub Actions job ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.843.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = "'1.768.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'328.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = "'0.4494"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.94%  '
$ws.Range("D8").Value = "'0.3560"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").Value = "'0.07460"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = "'42.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = "'21.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = "'6.034"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D16").Value = "'1.770.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").Value = "'93.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = "'0.06432"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = "'17.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.45%  '
$ws.Range("D22").Value = "'5.782"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").Value = "'27.900.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").Value = "'11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("D25").Value = "'2.119"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").Value = "'162.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").Value = "'20.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").Value = "'1.976.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").Value = "'2.167"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.26%  '
$ws.Range("D30").Value = "'125.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("E31").Value = '  +4.59%  '
$ws.Range("D32").Value = "'0.09188"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("D34").Value = "'3.650"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = "'11.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").Value = "'0.02293"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").Value = "'0.06137"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.92%  '
$ws.Range("D38").Value = "'0.2103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.99%  '
$ws.Range("D39").Value = "'0.6338"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").Value = "'4.968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").Value = "'1.185"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("E42").Value = '  +1.29%  '
$ws.Range("D43").Value = "'7.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").Value = "'13.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").Value = "'3.741"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").Value = "'0.5889"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").Value = "'122.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = "'1.957"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = "'0.06906"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.73%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").Value = "'73.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.96%  '
